$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 0.077
$ws.Cells.Item(2, 5).Value = 0.095
$ws.Cells.Item(2, 6).Value = 0.08

$ws.Cells.Item(3, 4).Value = 0.042
$ws.Cells.Item(3, 5).Value = 0.14
$ws.Cells.Item(3, 6).Value = 0.062

$ws.Cells.Item(4, 4).Value = 0.007
$ws.Cells.Item(4, 6).Value = 0.007

$ws.Cells.Item(5, 4).Value = 0.442
$ws.Cells.Item(5, 5).Value = 0.468
$ws.Cells.Item(5, 6).Value = 0.449

$ws.Cells.Item(6, 4).Value = 0.068
$ws.Cells.Item(6, 5).Value = 0.088
$ws.Cells.Item(6, 6).Value = 0.072
$ws.Cells.Item(6, 8).Value = 532

$ws.Cells.Item(7, 4).Value = 0.216
$ws.Cells.Item(7, 5).Value = 0.219
$ws.Cells.Item(7, 6).Value = 0.216

$ws.Cells.Item(8, 4).Value = 0.601
$ws.Cells.Item(8, 5).Value = 0.609
$ws.Cells.Item(8, 6).Value = 0.592

$ws.Cells.Item(9, 4).Value = 0.299
$ws.Cells.Item(9, 5).Value = 0.326
$ws.Cells.Item(9, 6).Value = 0.309
$ws.Cells.Item(9, 8).Value = 1169

$ws.Cells.Item(10, 4).Value = 0.295
$ws.Cells.Item(10, 5).Value = 0.327
$ws.Cells.Item(10, 6).Value = 0.306
$ws.Cells.Item(10, 8).Value = 1169

$ws.Cells.Item(11, 4).Value = 0.837
$ws.Cells.Item(11, 5).Value = 1.679
$ws.Cells.Item(11, 6).Value = 0.994

$ws.Cells.Item(12, 4).Value = 0.313
$ws.Cells.Item(12, 5).Value = 0.329
$ws.Cells.Item(12, 6).Value = 0.314

$ws.Cells.Item(13, 4).Value = 0.263
$ws.Cells.Item(13, 5).Value = 0.345
$ws.Cells.Item(13, 6).Value = 0.282

$ws.Cells.Item(14, 4).Value = 0.21
$ws.Cells.Item(14, 5).Value = 0.232
$ws.Cells.Item(14, 6).Value = 0.214

$ws.Cells.Item(15, 4).Value = 0.757
$ws.Cells.Item(15, 5).Value = 0.822
$ws.Cells.Item(15, 6).Value = 0.771

$ws.Cells.Item(16, 4).Value = 0.697
$ws.Cells.Item(16, 5).Value = 0.754
$ws.Cells.Item(16, 6).Value = 0.702

$ws.Cells.Item(17, 4).Value = 0.371
$ws.Cells.Item(17, 5).Value = 0.399
$ws.Cells.Item(17, 6).Value = 0.378

$ws.Cells.Item(18, 4).Value = 0.575
$ws.Cells.Item(18, 5).Value = 0.671
$ws.Cells.Item(18, 6).Value = 0.592

$ws.Cells.Item(19, 4).Value = 0.308
$ws.Cells.Item(19, 5).Value = 0.341
$ws.Cells.Item(19, 6).Value = 0.321
$ws.Cells.Item(19, 8).Value = 1169

$ws.Cells.Item(20, 4).Value = 0.307
$ws.Cells.Item(20, 5).Value = 0.336
$ws.Cells.Item(20, 6).Value = 0.316
$ws.Cells.Item(20, 8).Value = 1169

$ws.Cells.Item(21, 4).Value = 0.668
$ws.Cells.Item(21, 5).Value = 1.69
$ws.Cells.Item(21, 6).Value = 0.853

$ws.Cells.Item(22, 4).Value = 0.043
$ws.Cells.Item(22, 5).Value = 0.044
$ws.Cells.Item(22, 6).Value = 0.043

$ws.Cells.Item(23, 5).Value = 0.036
$ws.Cells.Item(23, 6).Value = 0.034

$ws.Cells.Item(24, 5).Value = 0.007
$ws.Cells.Item(24, 6).Value = 0.007

$ws.Cells.Item(25, 4).Value = 0.619
$ws.Cells.Item(25, 5).Value = 0.634
$ws.Cells.Item(25, 6).Value = 0.613

$ws.Cells.Item(26, 5).Value = 0.047
$ws.Cells.Item(26, 6).Value = 0.043
$ws.Cells.Item(26, 8).Value = 336

$ws.Cells.Item(27, 4).Value = 0.286
$ws.Cells.Item(27, 5).Value = 0.343
$ws.Cells.Item(27, 6).Value = 0.299

$ws.Cells.Item(28, 4).Value = 0.316
$ws.Cells.Item(28, 5).Value = 0.344
$ws.Cells.Item(28, 6).Value = 0.325

$ws.Cells.Item(29, 4).Value = 0.033
$ws.Cells.Item(29, 5).Value = 0.036
$ws.Cells.Item(29, 6).Value = 0.034

$ws.Cells.Item(30, 4).Value = 0.007
$ws.Cells.Item(30, 5).Value = 0.009
$ws.Cells.Item(30, 6).Value = 0.007

$ws.Cells.Item(31, 4).Value = 0.241
$ws.Cells.Item(31, 5).Value = 0.272
$ws.Cells.Item(31, 6).Value = 0.246

$ws.Cells.Item(32, 4).Value = 0.115
$ws.Cells.Item(32, 5).Value = 0.117
$ws.Cells.Item(32, 6).Value = 0.115

$ws.Cells.Item(33, 4).Value = 0.212
$ws.Cells.Item(33, 5).Value = 0.235
$ws.Cells.Item(33, 6).Value = 0.216

$ws.Cells.Item(34, 4).Value = 0.107
$ws.Cells.Item(34, 5).Value = 0.107
$ws.Cells.Item(34, 6).Value = 0.107

$ws.Cells.Item(35, 4).Value = 0.629
$ws.Cells.Item(35, 5).Value = 0.66
$ws.Cells.Item(35, 6).Value = 0.624

$ws.Cells.Item(36, 4).Value = 0.124
$ws.Cells.Item(36, 5).Value = 0.145
$ws.Cells.Item(36, 6).Value = 0.129

$ws.Cells.Item(37, 4).Value = 0.295
$ws.Cells.Item(37, 5).Value = 0.359
$ws.Cells.Item(37, 6).Value = 0.307

$ws.Cells.Item(38, 4).Value = 0.324
$ws.Cells.Item(38, 5).Value = 0.35
$ws.Cells.Item(38, 6).Value = 0.331

$ws.Cells.Item(39, 5).Value = 0.196
$ws.Cells.Item(39, 6).Value = 0.194

$ws.Cells.Item(40, 4).Value = 0.194
$ws.Cells.Item(40, 5).Value = 0.196
$ws.Cells.Item(40, 6).Value = 0.195

$ws.Cells.Item(41, 4).Value = 0.241
$ws.Cells.Item(41, 5).Value = 0.262
$ws.Cells.Item(41, 6).Value = 0.245

Write-Output "Applied 121 cell updates"
